$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'67.900.95"
$ws.Range("E2").Value = "  +0.93%  "

# Row 3
$ws.Range("D3").Value = "'2.513.21"
$ws.Range("E3").Value = "  -1.28%  "

# Row 5
$ws.Range("D5").Value = "'593.44"
$ws.Range("E5").Value = "  +0.42%  "

# Row 6
$ws.Range("D6").Value = "'173.83"
$ws.Range("E6").Value = "  -0.33%  "

# Row 7
$ws.Range("E7").Value = "  -0.01%  "

# Row 8
$ws.Range("D8").Value = "'0.528"
$ws.Range("E8").Value = "  -0.34%  "

# Row 9
$ws.Range("D9").Value = "'2.515.26"
$ws.Range("E9").Value = "  -1.23%  "

# Row 10
$ws.Range("D10").Value = "'0.139"
$ws.Range("E10").Value = "  -0.42%  "

# Row 12
$ws.Range("D12").Value = "'5.10"
$ws.Range("E12").Value = "  -1.19%  "

# Row 13
$ws.Range("D13").Value = "'0.343"
$ws.Range("E13").Value = "  -2.53%  "

# Row 14
$ws.Range("D14").Value = "'26.40"
$ws.Range("E14").Value = "  -2.25%  "

# Row 15
$ws.Range("D15").Value = "'2.973.96"
$ws.Range("E15").Value = "  -1.26%  "

# Row 16
$ws.Range("D16").Value = "'0.0000177"
$ws.Range("E16").Value = "  -0.50%  "

# Row 17
$ws.Range("D17").Value = "'67.767.28"
$ws.Range("E17").Value = "  +1.01%  "

# Row 18
$ws.Range("D18").Value = "'2.492.15"
$ws.Range("E18").Value = "  -2.40%  "

# Row 19
$ws.Range("D19").Value = "'11.84"
$ws.Range("E19").Value = "  +3.78%  "

# Row 20
$ws.Range("D20").Value = "'7.97"
$ws.Range("E20").Value = "  -1.43%  "

# Row 21
$ws.Range("D21").Value = "'363.43"
$ws.Range("E21").Value = "  +2.18%  "

# Row 22
$ws.Range("D22").Value = "'4.13"
$ws.Range("E22").Value = "  -2.01%  "

# Row 23
$ws.Range("D23").Value = "'4.61"
$ws.Range("E23").Value = "  -1.30%  "

# Row 24
$ws.Range("D24").Value = "'71.29"
$ws.Range("E24").Value = "  +1.75%  "

# Row 25
$ws.Range("E25").Value = "  -0.07%  "

# Row 26
$ws.Range("D26").Value = "'1.91"
$ws.Range("E26").Value = "  -4.41%  "

# Row 27
$ws.Range("D27").Value = "'9.94"
$ws.Range("E27").Value = "  -0.93%  "

# Row 28
$ws.Range("E28").Value = "  +0.05%  "

# Row 29
$ws.Range("D29").Value = "'2.645.62"
$ws.Range("E29").Value = "  -1.47%  "

# Row 30
$ws.Range("D30").Value = "'0.0₃0975"
$ws.Range("E30").Value = "  -2.28%  "

# Row 31
$ws.Range("D31").Value = "'8.30"
$ws.Range("E31").Value = "  +0.39%  "

# Row 32
$ws.Range("D32").Value = "'530.61"
$ws.Range("E32").Value = "  -0.90%  "

# Row 33
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").Value = "'1.32"
$ws.Range("E33").Value = "  -2.31%  "

# Row 34
$ws.Range("B34").Value = "PancakeSwap"
$ws.Range("C34").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D34").Value = "'1.88"
$ws.Range("E34").Value = "  +1.21%  "

# Row 35
$ws.Range("D35").Value = "'0.128"
$ws.Range("E35").Value = "  -2.80%  "

# Row 36
$ws.Range("E36").Value = "  +0.03%  "

# Row 37
$ws.Range("D37").Value = "'158.50"
$ws.Range("E37").Value = "  +0.62%  "

# Row 38
$ws.Range("D38").Value = "'1.44"
$ws.Range("E38").Value = "  -2.94%  "

# Row 39
$ws.Range("D39").Value = "'18.68"
$ws.Range("E39").Value = "  -0.67%  "

# Row 40
$ws.Range("E40").Value = "  +1.18%  "

# Row 41
$ws.Range("D41").Value = "'1.79"
$ws.Range("E41").Value = "  -1.14%  "

# Row 42
$ws.Range("B42").Value = "PolygonEcosystemToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D42").Value = "'0.350"
$ws.Range("E42").Value = "  -2.11%  "

# Row 43
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D43").Value = "'5.13"
$ws.Range("E43").Value = "  -1.55%  "

# Row 44
$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D44").Value = "'0.999"
$ws.Range("E44").Value = "  -0.15%  "

# Row 45
$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").Value = "'2.50"
$ws.Range("E45").Value = "  -2.25%  "

# Row 46
$ws.Range("D46").Value = "'146.73"
$ws.Range("E46").Value = "  -3.01%  "

# Row 47
$ws.Range("B47").Value = "Filecoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D47").Value = "'3.70"
$ws.Range("E47").Value = "  -0.56%  "

# Row 48
$ws.Range("D48").Value = "'0.551"
$ws.Range("E48").Value = "  -2.47%  "

# Row 49
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "'0.0₆0276"
$ws.Range("E49").Value = "  -1.86%  "

# Row 50
$ws.Range("D50").Value = "'1.72"
$ws.Range("E50").Value = "  -0.30%  "

# Row 51
$ws.Range("E51").Value = "  -1.13%  "
